$wb = $excel.ActiveWorkbook

# Sheet "A1"
$ws = $wb.Worksheets.Item("A1")
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 0
$ws.Range("D1").Value = 18.4
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = [double]"6.123233995736766e-17"
$ws.Range("C2").Value = -1
$ws.Range("D2").Value = 0

# Sheet "A2"
$ws = $wb.Worksheets.Item("A2")
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 0
$ws.Range("D1").Value = 149
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1
$ws.Range("D2").Value = 0

# Sheet "A3"
$ws = $wb.Worksheets.Item("A3")
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 0
$ws.Range("D1").Value = 120.3
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1
$ws.Range("D2").Value = 0

# Sheet "A4"
$ws = $wb.Worksheets.Item("A4")
$ws.Range("A1").Value = [double]"6.123233995736766e-17"
$ws.Range("B1").Value = [double]"6.123233995736766e-17"
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = [double]"5.37619944825688e-15"
$ws.Range("A2").Value = -1
$ws.Range("B2").Value = [double]"3.749399456654644e-33"
$ws.Range("C2").Value = [double]"6.123233995736766e-17"
$ws.Range("D2").Value = -87.8

# Sheet "A_total"
$ws = $wb.Worksheets.Item("A_total")
$ws.Range("A1").Value = [double]"6.123233995736766e-17"
$ws.Range("B1").Value = [double]"6.123233995736766e-17"
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 297.7
$ws.Range("A2").Value = [double]"-6.123233995736766e-17"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = [double]"-6.123233995736766e-17"
$ws.Range("D2").Value = [double]"-7.396866666850013e-15"
$ws.Range("A3").Value = -1
$ws.Range("B3").Value = [double]"-6.123233995736766e-17"
$ws.Range("C3").Value = [double]"6.123233995736766e-17"
$ws.Range("D3").Value = 0
